# Apply the 27-Jul-2024 cryptos list refresh to Sheet1.
# Column D holds price text that sometimes LOOKS like a plain number
# (e.g. "1.01"); assigning such text via .Value would make Excel
# auto-convert it to a real number, so those writes are given a
# leading apostrophe to force a text entry, matching the source data's
# inlineStr cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.088.87'
$ws.Range("E2").Value = '  +1.25%  '

$ws.Range("D3").Value = '3.267.83'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''586.27'
$ws.Range("E5").Value = '  +1.76%  '

$ws.Range("D6").Value = '''184.68'
$ws.Range("E6").Value = '  +3.56%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '''0.599'
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +4.53%  '

$ws.Range("D10").Value = '''6.73'
$ws.Range("E10").Value = '  +0.22%  '

$ws.Range("D11").Value = '''0.416'
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").Value = '3.838.49'
$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("E13").Value = '  +0.31%  '

$ws.Range("D14").Value = '''28.59'

$ws.Range("D15").Value = '68.092.41'

$ws.Range("E16").Value = '  +2.22%  '

$ws.Range("D17").Value = '3.267.99'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").Value = '''13.61'
$ws.Range("E19").Value = '  +2.35%  '

$ws.Range("D20").Value = '''382.25'
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("D21").Value = '''7.73'
$ws.Range("E21").Value = '  +1.75%  '

$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("E23").Value = '  +0.68%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '''0.514'
$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000121'
$ws.Range("E25").Value = '  +2.59%  '

$ws.Range("E26").Value = '  +6.72%  '

$ws.Range("D27").Value = '''9.77'
$ws.Range("E27").Value = '  -1.24%  '

$ws.Range("D28").Value = '''1.01'
$ws.Range("E28").Value = '  +0.43%  '

$ws.Range("D29").Value = '''5.81'
$ws.Range("E29").Value = '  +3.76%  '

$ws.Range("D30").Value = '''1.99'
$ws.Range("E30").Value = '  +1.22%  '

$ws.Range("D31").Value = '''22.92'
$ws.Range("E31").Value = '  +1.80%  '

$ws.Range("D32").Value = '''7.19'
$ws.Range("E32").Value = '  +5.75%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("E34").Value = '  +0.21%  '

$ws.Range("D35").Value = '''1.54'
$ws.Range("E35").Value = '  +3.09%  '

$ws.Range("D36").Value = '''163.36'
$ws.Range("E36").Value = '  +0.93%  '

$ws.Range("E37").Value = '  +0.76%  '

$ws.Range("D38").Value = '''0.838'
$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("E39").Value = '  -1.52%  '

$ws.Range("D40").Value = '''26.58'
$ws.Range("E40").Value = '  -0.83%  '

$ws.Range("D41").Value = '''2.64'
$ws.Range("E41").Value = '  +1.39%  '

$ws.Range("D42").Value = '''4.59'
$ws.Range("E42").Value = '  +4.64%  '

$ws.Range("D43").Value = '''25.59'
$ws.Range("E43").Value = '  -0.53%  '

$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = '''0.0690'
$ws.Range("E44").Value = '  +2.87%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''41.22'
$ws.Range("E45").Value = '  +1.93%  '

$ws.Range("D46").Value = '2.634.48'
$ws.Range("E46").Value = '  -4.61%  '

$ws.Range("D47").Value = '''341.74'
$ws.Range("E47").Value = '  -3.43%  '

$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("D49").Value = '''32.19'
$ws.Range("E49").Value = '  +4.89%  '

$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '  +1.24%  '

$ws.Range("E51").Value = '  -0.13%  '
